# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to reflect the freshly scraped data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 3..6 hold the F column counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2032
$wsExhibit.Range("F4").Value = 848
$wsExhibit.Range("F5").Value = 1099
$wsExhibit.Range("F6").Value = 348

# Sheet "全部类型" (All types) - same underlying rows, but shifted down by the
# extra "演出" rows merged in, so they land on F3 and F6..F8
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2032
$wsAll.Range("F6").Value = 848
$wsAll.Range("F7").Value = 1099
$wsAll.Range("F8").Value = 348
